# week 16 Meeting.pptx edit script
# - refresh the "fixed" datetime1 date-placeholder text (master + all layouts)
#   from 2/7/2022 -> 4/8/2022
# - update the title-slide subtitle text to the new meeting number/date

$p = $ppt.ActivePresentation
$cr = [char]13

$oldDate = "2/7/2022"
$newDate = "4/8/2022"

# --- Slide master: date placeholder text -------------------------------
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide layouts: date placeholder text -------------------------------
$layouts = $m.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 1: subtitle "Week 14 Meeting, 31/1/2021" -> "Week 15 Meeting, 6/2/2021"
$oldSubtitle = "Week 14 Meeting, 31/1/2021" + $cr + $cr
$newSubtitle = "Week 15 Meeting, 6/2/2021" + $cr + $cr

$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldSubtitle) {
            # Replace via a neutral placeholder first so the engine performs a
            # single clean run substitution instead of diff-splitting the run
            # against the old characters.
            $sh.TextFrame.TextRange.Text = "x" + $cr + $cr
            $sh.TextFrame.TextRange.Text = $newSubtitle
        }
    }
}
